$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$data = @(
    @("ECs", "Efna5", "Epha4", "ECs", 2, 0.6666666666666666, 0.230855, 0.692565, 0.06377305075821572, 0.06377305075821572, 2, 0.6666666666666666, 6.708176333333333, 20.124529, 0.4356329228871633, 0.4356329228871633, 1.548616047431667, 13.937544426885, 0.02778164050323294, 0.02778164050323294),
    @("ECs", "Efna5", "Epha4", "FAPs", 2, 0.6666666666666666, 0.230855, 0.692565, 0.06377305075821572, 0.06377305075821572, 3, 1, 6.789877333333333, 20.369632, 0.4409386339573907, 0.4409386339573907, 1.567477131786667, 14.10729418608, 0.02812000188462298, 0.02812000188462298),
    @("ECs", "Efna5", "Epha4", "sCs", 2, 0.6666666666666666, 0.230855, 0.692565, 0.06377305075821572, 0.06377305075821572, 3, 1, 1.900636333333334, 5.701909000000001, 0.1234284431554459, 0.1234284431554459, 0.4387714007316668, 3.948942606585, 0.007871408370359797, 0.007871408370359797),
    @("FAPs", "Efna5", "Epha4", "ECs", 3, 1, 2.900731333333333, 8.702194, 0.8013189515350044, 0.8013189515350045, 2, 0.6666666666666666, 6.708176333333333, 20.124529, 0.4356329228871633, 0.4356329228871633, 19.45861727962511, 175.127555516626, 0.3490809170220712, 0.3490809170220712),
    @("FAPs", "Efna5", "Epha4", "FAPs", 3, 1, 2.900731333333333, 8.702194, 0.8013189515350044, 0.8013189515350045, 3, 1, 6.789877333333333, 20.369632, 0.4409386339573907, 0.4409386339573907, 19.69560993028978, 177.260489372608, 0.3533324838540134, 0.3533324838540135),
    @("FAPs", "Efna5", "Epha4", "sCs", 3, 1, 2.900731333333333, 8.702194, 0.8013189515350044, 0.8013189515350045, 3, 1, 1.900636333333334, 5.701909000000001, 0.1234284431554459, 0.1234284431554459, 5.513235365371779, 49.61911828834601, 0.09890555065891982, 0.09890555065891983),
    @("sCs", "Efna5", "Epha4", "ECs", 3, 1, 0.4883596666666667, 1.465079, 0.1349079977067798, 0.1349079977067798, 2, 0.6666666666666666, 6.708176333333333, 20.124529, 0.4356329228871633, 0.4356329228871633, 3.276002758087889, 29.484024822791, 0.05877036536185921, 0.05877036536185921),
    @("sCs", "Efna5", "Epha4", "FAPs", 3, 1, 0.4883596666666667, 1.465079, 0.1349079977067798, 0.1349079977067798, 3, 1, 6.789877333333333, 20.369632, 0.4409386339573907, 0.4409386339573907, 3.315902231214222, 29.843120080928, 0.05948614821875429, 0.05948614821875429),
    @("sCs", "Efna5", "Epha4", "sCs", 3, 1, 0.4883596666666667, 1.465079, 0.1349079977067798, 0.1349079977067798, 3, 1, 1.900636333333334, 5.701909000000001, 0.1234284431554459, 0.1234284431554459, 0.9281941262012224, 8.353747135811002, 0.0166514841261663, 0.0166514841261663),
)

$startRow = 2
for ($i = 0; $i -lt $data.Count; $i++) {
    $row = $data[$i]
    for ($j = 0; $j -lt $row.Count; $j++) {
        $ws.Cells.Item($startRow + $i, $j + 1).Value = $row[$j]
    }
}
